# Appends a new data row (row 9) to the "Artfynd" export sheet, matching a
# fresh observation record (Goodyera repens / Knärot) that was added to the
# source data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- plain numeric cells -----------------------------------------------
$ws.Range("A9").Value  = 112486937
$ws.Range("B9").Value  = 96720
$ws.Range("E9").Value  = 220787
$ws.Range("Q9").Value  = 428634
$ws.Range("R9").Value  = 6274071
$ws.Range("S9").Value  = 18

# --- plain text cells -----------------------------------------------------
$ws.Range("C9").Value  = "Ovaliderad"
$ws.Range("D9").Value  = "VU"
$ws.Range("F9").Value  = "Knärot"
$ws.Range("G9").Value  = "Goodyera repens"
$ws.Range("H9").Value  = "(L.) R. Br."
$ws.Range("J9").Value  = "stjälkar/strån/skott"
$ws.Range("K9").Value  = "överblommad"
$ws.Range("P9").Value  = "Hyltåkra, Sm"
$ws.Range("T9").Value  = "Kronoberg"
$ws.Range("U9").Value  = "Ljungby"
$ws.Range("V9").Value  = "Småland"
$ws.Range("W9").Value  = "Hamneda"
$ws.Range("X9").Value  = "G-Lju-1389"
$ws.Range("AC9").Value = "1 blomställning. Rullstensås. Växer på åsens norra slänt. Olikåldrig barrskog. Tallarna säkert 150 år gamla, granarna något yngre (50 - 100 år)."
$ws.Range("AW9").Value = "Krister Wahlström"
$ws.Range("AX9").Value = "Krister Wahlström"
$ws.Range("AY9").Value = "Floraväkteri Sverige"

# --- text cells whose content would otherwise auto-convert (numeric- or
#     date-looking strings) - a leading apostrophe forces text, same as a
#     user typing it in Excel. ------------------------------------------
$ws.Range("I9").Value  = "'16"
$ws.Range("Y9").Value  = "'2023-10-02"
$ws.Range("AA9").Value = "'2023-10-02"

# --- boolean cells ----------------------------------------------------
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
